$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Notas" column (J), formats copied from column I so styles/xfs are reused ---
$ws.Range("I4:I11").Copy()
$ws.Range("J4:J11").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("J4").Value = "Notas"
$ws.Range("J5").Value = "N/A"
$ws.Range("J6").Value = "N/A"

# Approximate column J width to match column I (closest reproducible value)
$ws.Columns.Item(10).ColumnWidth = 40.333333333333336

# --- New data row for 05/05/2025 (Lunes 05 Mayo) ---
# Copy formatting from row 6 (the previous data row) onto row 7 first, then set values
$ws.Range("D6:I6").Copy()
$ws.Range("D7:I7").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("D7").Value = 45782
$ws.Range("E7").Value = 585
$ws.Range("F7").Value = 202
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 300
$ws.Range("J7").Value = "Se subieron 300 imagenes a dataset "

# Update the currently-selected cell shown in the workbook view
$ws.Range("K6").Select()
